# Notes_ApresRat.xlsx - "ajout de fichier de conversion Excel en xml"
#
# The sheet header row had its LastName/FirstName columns swapped (B1/C1),
# and the student CNE numbers in column A (rows 3-11) were de-duplicated
# into a proper incrementing sequence. The active selection also moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Header row: swap B1 (LastName) and C1 (FirstName) ---
$colB = $ws.Range("B1").Value2
$colC = $ws.Range("C1").Value2
$ws.Range("B1").Value = $colC
$ws.Range("C1").Value = $colB

# --- Column A (CNE): rows 3..11 go from the duplicated 19000041 to a
#     proper incrementing sequence 19000042..19000050. Row 2 is untouched. ---
for ($r = 3; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = 19000039 + $r
}

# --- Default row height / font metrics used by the sheet ---
$ws.StandardHeight = 14.5

# --- Active selection moves from L9 to F7 ---
$ws.Range("F7").Select() | Out-Null
